$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Numeric-looking values (Price / Volume %) must stay as text, matching the original inlineStr cells.
# Temporarily apply a text number format so Excel does not auto-convert them to numbers,
# then restore the original (default/General) formatting afterwards.
$numericCells = @("D2","E2","D3","E3","D4","E4","E5","D6","E6","D7","E7","D8","E8","D9","E9","E10","D11","E11","D12","E12","D13","E13","D14","E14","D15","E15","D16","E16","D17","E17","D18","E18","D19","E19","D20","E20","E21","D22","E22","D23","E23","D24","E24","D25","E25","E26","D38","E38","D39","E39","D40","E40","E41","D42","E42","D43","E43","D44","E44","D45","E45","D46","E46","D47","E47","E48","E49","E50")
foreach ($addr in $numericCells) {
    $ws.Range($addr).NumberFormat = "@"
}

$ws.Range("D2").Value = "326.57"
$ws.Range("E2").Value = "-2.50%"
$ws.Range("D3").Value = "44.51"
$ws.Range("E3").Value = "1.42%"
$ws.Range("D4").Value = "5.604"
$ws.Range("E4").Value = "-2.16%"
$ws.Range("E5").Value = "-3.22%"
$ws.Range("D6").Value = "4.300"
$ws.Range("E6").Value = "-4.79%"
$ws.Range("D7").Value = "1.897"
$ws.Range("E7").Value = "-3.22%"
$ws.Range("D8").Value = "2.650"
$ws.Range("E8").Value = "-7.71%"
$ws.Range("D9").Value = "0.9470"
$ws.Range("E9").Value = "0.33%"
$ws.Range("E10").Value = "-6.98%"
$ws.Range("D11").Value = "0.1846"
$ws.Range("E11").Value = "-6.96%"
$ws.Range("D12").Value = "0.09908"
$ws.Range("E12").Value = "-7.78%"
$ws.Range("D13").Value = "0.04289"
$ws.Range("E13").Value = "-5.50%"
$ws.Range("D14").Value = "0.1068"
$ws.Range("E14").Value = "0.11%"
$ws.Range("D15").Value = "0.001286"
$ws.Range("E15").Value = "-0.51%"
$ws.Range("D16").Value = "0.04220"
$ws.Range("E16").Value = "-4.59%"
$ws.Range("D17").Value = "0.005939"
$ws.Range("E17").Value = "0.46%"
$ws.Range("D18").Value = "3.610"
$ws.Range("E18").Value = "3.24%"
$ws.Range("D19").Value = "0.3498"
$ws.Range("E19").Value = "-0.20%"
$ws.Range("D20").Value = "8.445"
$ws.Range("E20").Value = "-2.52%"
$ws.Range("E21").Value = "1.38%"
$ws.Range("D22").Value = "0.2653"
$ws.Range("E22").Value = "-1.47%"
$ws.Range("D23").Value = "0.001249"
$ws.Range("E23").Value = "-0.56%"
$ws.Range("D24").Value = "0.004528"
$ws.Range("E24").Value = "4.20%"
$ws.Range("D25").Value = "0.0001262"
$ws.Range("E25").Value = "0.06%"
$ws.Range("E26").Value = "0.09%"
$ws.Range("D38").Value = "0.02633"
$ws.Range("E38").Value = "-6.33%"
$ws.Range("D39").Value = "0.05481"
$ws.Range("E39").Value = "-8.75%"
$ws.Range("D40").Value = "0.007624"
$ws.Range("E40").Value = "-3.66%"
$ws.Range("E41").Value = "-1.97%"
$ws.Range("D42").Value = "0.007358"
$ws.Range("E42").Value = "-17.98%"
$ws.Range("D43").Value = "0.002051"
$ws.Range("E43").Value = "-5.54%"
$ws.Range("D44").Value = "0.008833"
$ws.Range("E44").Value = "-12.72%"
$ws.Range("D45").Value = "0.00006927"
$ws.Range("E45").Value = "-1.17%"
$ws.Range("D46").Value = "0.00000000752"
$ws.Range("E46").Value = "0.09%"
$ws.Range("D47").Value = "0.003706"
$ws.Range("E47").Value = "16.30%"
$ws.Range("E48").Value = "0.09%"
$ws.Range("E49").Value = "0.09%"
$ws.Range("E50").Value = "0.09%"

foreach ($addr in $numericCells) {
    $ws.Range($addr).NumberFormat = "General"
    $ws.Range($addr).Style = "Normal"
}

# Plain text cells (coin name / link) - Excel keeps these as text natively.
$ws.Range("B6").Value = "GateToken"
$ws.Range("C6").Value = "https://coinranking.com/coin/t7m8DZVyMsAu+gatetoken-gt"
$ws.Range("B7").Value = "FTXToken"
$ws.Range("C7").Value = "https://coinranking.com/coin/NfeOYfNcl+ftxtoken-ftt"
